$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(3)
$shape.TextFrame.TextRange.Text = "contd.."
